$wb = $excel.ActiveWorkbook

# --- Meetings sheet: mark additional "ü" (hours-worked) cells for E and F
# columns on every task row (4-9). This raises each team lead's logged
# meeting hours total (column B) which flows into the Management Summary
# sheet via `Meetings!B4*100` etc. ---
$meetings = $wb.Worksheets.Item("Meetings")
$meetings.Range("E4").Value = "ü"
$meetings.Range("F4").Value = "ü"
$meetings.Range("E5").Value = "ü"
$meetings.Range("F5").Value = "ü"
$meetings.Range("E6").Value = "ü"
$meetings.Range("F6").Value = "ü"
$meetings.Range("E7").Value = "ü"
$meetings.Range("F7").Value = "ü"
$meetings.Range("E8").Value = "ü"
$meetings.Range("F8").Value = "ü"
$meetings.Range("E9").Value = "ü"
$meetings.Range("F9").Value = "ü"

# --- SA (Systems Analysis) sheet: fill in the budgeted/actual hour counts
# for each team lead's group (columns C and D of the first detail row of
# each group). These feed the Management Summary sheet via `SA!C5`,
# `SA!D5`, `SA!C9`, `SA!D9`, etc. ---
$sa = $wb.Worksheets.Item("SA")
$sa.Range("C2").Value = 2
$sa.Range("D2").Value = 2
$sa.Range("C6").Value = 1
$sa.Range("D6").Value = 1.5
$sa.Range("C10").Value = 2
$sa.Range("D10").Value = 3
$sa.Range("C16").Value = 2
$sa.Range("D16").Value = 1.5
$sa.Range("C20").Value = 2
$sa.Range("D20").Value = 2
$sa.Range("C24").Value = 2
$sa.Range("D24").Value = 2

# --- Restore cursor/selection on SA before moving away from it, matching
# the last-used cell recorded for that sheet. ---
[void]$sa.Activate()
[void]$sa.Range("D7").Select()

# --- Finally land on the Meetings sheet/cell, which becomes the active
# (saved) tab, mirroring the workbook's recorded view state. ---
[void]$meetings.Activate()
[void]$meetings.Range("F9").Select()
